$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> ECs)
$ws.Range("G2").Value = 0.1038603333333333
$ws.Range("H2").Value = 0.311581
$ws.Range("M2").Value = 0.165747
$ws.Range("N2").Value = 0.497241
$ws.Range("O2").Value = 0.008095785894995438
$ws.Range("P2").Value = 0.00809578589499544
$ws.Range("Q2").Value = 0.017214538669
$ws.Range("R2").Value = 0.154930848021
$ws.Range("S2").Value = 0.008095785894995438
$ws.Range("T2").Value = 0.00809578589499544

# Row 3 (MuSCs -> FAPs)
$ws.Range("G3").Value = 0.1038603333333333
$ws.Range("H3").Value = 0.311581
$ws.Range("O3").Value = 0.7079722685862583
$ws.Range("P3").Value = 0.7079722685862583
$ws.Range("Q3").Value = 1.505402459036334
$ws.Range("R3").Value = 13.548622131327
$ws.Range("S3").Value = 0.7079722685862583
$ws.Range("T3").Value = 0.7079722685862583

# Row 4 (MuSCs -> MuSCs)
$ws.Range("G4").Value = 0.1038603333333333
$ws.Range("H4").Value = 0.311581
$ws.Range("M4").Value = 5.642879333333333
$ws.Range("N4").Value = 16.928638
$ws.Range("O4").Value = 0.2756221404547972
$ws.Range("P4").Value = 0.2756221404547972
$ws.Range("Q4").Value = 0.5860713285197777
$ws.Range("R4").Value = 5.274641956678
$ws.Range("S4").Value = 0.2756221404547972
$ws.Range("T4").Value = 0.2756221404547972

# Row 5 (MuSCs -> Resolving-Mac)
$ws.Range("G5").Value = 0.1038603333333333
$ws.Range("H5").Value = 0.311581
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1701286666666667
$ws.Range("N5").Value = 0.510386
$ws.Range("O5").Value = 0.008309805063949155
$ws.Range("P5").Value = 0.008309805063949155
$ws.Range("Q5").Value = 0.01766962002955556
$ws.Range("R5").Value = 0.159026580266
$ws.Range("S5").Value = 0.008309805063949155
$ws.Range("T5").Value = 0.008309805063949155
